$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values for the rows that were re-pulled/recalculated.
$ws.Range("F5").Value = -3
$ws.Range("F11").Value = -2
$ws.Range("F14").Value = 6
$ws.Range("F15").Value = 1
$ws.Range("F23").Value = -3
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = -4
$ws.Range("F27").Value = 7
$ws.Range("F28").Value = -6
